$wb = $excel.ActiveWorkbook

# --- Overview sheet: status message changes from "Ready for handoff" to
#     "Handed back: in sync with en-US" for both tracked files ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

# Widen the two status columns to fit the longer text
$overview.Columns.Item(5).ColumnWidth = 29.9777050018311
$overview.Columns.Item(6).ColumnWidth = 29.9777050018311

# --- zh-cn sheet: refresh handback datetime, clear stale error detail ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("K2").Value = "2016-10-18 12:19:39"
$zhcn.Range("K3").Value = "2016-10-18 12:19:39"
$zhcn.Range("P2").Value = ""
$zhcn.Range("P3").Value = ""

$zhcn.Columns.Item(3).ColumnWidth = 29.9777050018311
$zhcn.Columns.Item(16).ColumnWidth = 13.7470531463623

# --- de-de sheet: refresh handback datetime, clear stale error detail ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("K2").Value = "2016-10-18 12:19:56"
$dede.Range("K3").Value = "2016-10-18 12:19:56"
$dede.Range("P2").Value = ""
$dede.Range("P3").Value = ""

$dede.Columns.Item(3).ColumnWidth = 29.9777050018311
$dede.Columns.Item(16).ColumnWidth = 13.7470531463623
